$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Absent" (column H) report values.
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 0
